# Move the "Causality and the notion of ceteris paribus" slide
# (currently slide 10) so that it sits right after slide 5
# (i.e. becomes the new slide 6), pushing the slides that used to be
# 6-9 down by one position (they become 7-10). All other slides keep
# their relative order.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$s.MoveTo(6)
